$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("level2")

# Clear the AD column values (rows 1-12) while keeping the existing style.
for ($r = 1; $r -le 12; $r++) {
    $ws.Cells.Item($r, 30).ClearContents()
}

# Move the active selection on the "level2" sheet to AG4.
$ws.Activate()
$ws.Range("AG4").Select()
